$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values re-pulled from source data.
$updates = @{
    2  = -4
    4  = -1
    5  = -1
    7  = -1
    9  = 1
    10 = -3
    12 = -2
    13 = 6
    15 = 1
    20 = 0
    21 = -4
    23 = -2
    26 = -2
    28 = -3
    31 = -1
    32 = 2
    33 = 6
    35 = 0
    39 = -1
    42 = 2
    43 = 2
    44 = -4
    45 = 1
    46 = 5
    47 = 0
    49 = -4
    50 = -2
    51 = 4
    53 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
